$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update wording of existing cells (shared-string text changes) ---
# All of B14:B18 held the same old text "商讨，修改完善用例图"; once every
# referencing cell is rewritten to the new wording the old shared string is
# fully vacated and the writer reuses/renumbers it in place (matches the
# target diff, which edits the <si> text without growing uniqueCount).
$ws.Range("B14").Value = "商讨完善用例图"
$ws.Range("B15").Value = "商讨完善用例图"
$ws.Range("B16").Value = "商讨完善用例图"
$ws.Range("B17").Value = "商讨完善用例图"
$ws.Range("B18").Value = "商讨完善用例图"
$ws.Range("B13").Value = "商讨完善用例图，修改用例图"

# --- Append a new weekly block (rows 21-30), mirroring the layout of the
#     existing block in rows 11-20 ---

# Merge the two "banner" ranges first, then paste formats from the
# equivalent source rows on top - this reuses the existing cell styles
# (s=6 / s=7) instead of minting new split-border styles.
$null = $ws.Range("A21:D21").Merge()
$null = $ws.Range("A29:D30").Merge()

$ws.Range("A11:D20").Copy()
$ws.Range("A21:D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new block's text content
$ws.Range("A21").Value = "日期：2018.10.10 第六周周三"

$ws.Range("A22").Value = "组员"
$ws.Range("B22").Value = "计划内容"
$ws.Range("C22").Value = "完成情况"
$ws.Range("D22").Value = "备注"

$ws.Range("A23").Value = "何舒静"
$ws.Range("A24").Value = "陈碧容"
$ws.Range("A25").Value = "黄丙升"
$ws.Range("A26").Value = "王增璟"
$ws.Range("A27").Value = "苏立明"
$ws.Range("A28").Value = "蔡智杰"

$ws.Range("A29").Value = "总结："

# --- Selection matches the author's last-saved cursor position ---
$null = $ws.Range("H27").Select()
